$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clipe")
$lo = $ws.ListObjects.Item(1)

# Capture source (IdTransmissao) values before adding the new column
$srcRange = $lo.ListColumns.Item("IdTransmissao").DataBodyRange
$rowCount = $srcRange.Rows.Count
$values = @()
for ($i = 1; $i -le $rowCount; $i++) {
    $values += , ($srcRange.Cells.Item($i, 1).Value())
}

# Add the new table column and name it
$newCol = $lo.ListColumns.Add()
$ws.Range("H1").Value = "IdUsuarioEspectador"

# Fill the new column with the same values as IdTransmissao
$destRange = $newCol.DataBodyRange
for ($i = 1; $i -le $rowCount; $i++) {
    $destRange.Cells.Item($i, 1).Value = $values[$i - 1]
}

# Match the header formatting used elsewhere in the table: the new last
# header cell gets a light (background1/white) font colour, and the first
# header cell gets an underline.
$h1 = $ws.Range("H1")
$h1.Font.ThemeColor = 2
$h1.WrapText = $false

$a1 = $ws.Range("A1")
$a1.Font.Underline = $true
